# Insert a new worksheet "具有相當價值之財產" (Property of considerable value)
# between "存款" and "保險", carrying one row of jewelry-material data.
# All other sheets (保險, 債務, ...) keep their original content unchanged.

$wb = $excel.ActiveWorkbook

$deposits = $wb.Worksheets.Item("存款")

# Add the new sheet right after "存款" (i.e. right before "保險").
$newWs = $wb.Worksheets.Add($null, $deposits)
$newWs.Name = "具有相當價值之財產"

# ---- Header row (row 1) ----
$newWs.Range("B1").Value = "name"
$newWs.Range("C1").Value = "quantity"
$newWs.Range("D1").Value = "owner"
$newWs.Range("E1").Value = "total"
$newWs.Range("F1").Value = "property_category"
$newWs.Range("G1").Value = "category"
$newWs.Range("H1").Value = "date"
$newWs.Range("I1").Value = "legislator_name"
$newWs.Range("J1").Value = "legislator_id"
$newWs.Range("K1").Value = "source_file"
$newWs.Range("L1").Value = "index"

# ---- Data row (row 2) ----
$newWs.Range("A2").Value = 86
$newWs.Range("B2").Value = "珠寶材料"
$newWs.Range("C2").Value = 40
$newWs.Range("D2").Value = "黃素香"
$newWs.Range("E2").Value = "6000000(製作珠寶飾品之原材料（估計價值））"
$newWs.Range("F2").Value = "antique"
$newWs.Range("G2").Value = "normal"
# Use a leading apostrophe so the date-looking text stays plain text
# instead of being auto-converted into a date serial number.
$newWs.Range("H2").Value = "'2013-12-26"
$newWs.Range("I2").Value = "李桐豪"
$newWs.Range("J2").Value = 896
$newWs.Range("K2").Value = "tmp2e9d1"
$newWs.Range("L2").Value = 86

# ---- Formatting: mirror the style used by the other sheets ----
# Header row + the A2 index cell: bold, thin box border, centered/top aligned.
$headerRow = $newWs.Range("B1:L1")
$headerRow.Font.Bold = $true
$headerRow.Borders.LineStyle = 1
$headerRow.HorizontalAlignment = -4108
$headerRow.VerticalAlignment = -4160

$indexCell = $newWs.Range("A2")
$indexCell.Font.Bold = $true
$indexCell.Borders.LineStyle = 1
$indexCell.HorizontalAlignment = -4108
$indexCell.VerticalAlignment = -4160
